$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain text (matches original inline-string cells)
# so numeric-looking values like "0.110" or "580.21" keep their exact text form
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.774.33"
$ws.Range("E2").Value = "  +4.74%  "
$ws.Range("D3").Value = "2.978.57"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "580.21"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "152.29"
$ws.Range("E6").Value = "  +6.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "2.976.29"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("E14").Value = "  +6.25%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "64.728.06"
$ws.Range("E16").Value = "  +4.86%  "
$ws.Range("D17").Value = "3.473.94"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "6.88"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "2.978.69"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").Value = "445.59"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  +5.11%  "
$ws.Range("D24").Value = "80.80"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "10.79"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "2.37"
$ws.Range("E29").Value = "  +14.48%  "
$ws.Range("D30").Value = "7.70"
$ws.Range("E30").Value = "  +8.72%  "
$ws.Range("D31").Value = "0.0000104"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "2.57"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").Value = "26.57"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "0.979"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  +7.71%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "48.98"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("D41").Value = "43.45"
$ws.Range("E41").Value = "  +9.65%  "
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("D43").Value = "0.294"
$ws.Range("E43").Value = "  +8.99%  "
$ws.Range("D44").Value = "8.39"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "381.79"
$ws.Range("E45").Value = "  +12.54%  "
$ws.Range("D46").Value = "2.782.79"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "0.0348"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("D48").Value = "134.16"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "0.000217"
$ws.Range("E51").Value = "  +12.50%  "
